# Apply updated crypto price/volume figures (Sat Jul  8 08:17:02 UTC 2023 GitHub Actions run).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "30.220.19"
$ws.Range("E2").Value = "  +0.46%  "
$ws.Range("D3").Value = "1.862.75"
$ws.Range("E3").Value = "  +0.11%  "
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "236.21"
$ws.Range("E5").Value = "  +1.18%  "
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  -0.03%  "
$ws.Range("D7").Value = "0.4675"
$ws.Range("E7").Value = "  +0.49%  "
$ws.Range("D8").Value = "0.2851"
$ws.Range("E8").Value = "  +1.77%  "
$ws.Range("D9").Value = "0.06536"
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "22.39"
$ws.Range("E10").Value = "  +14.73%  "
$ws.Range("D11").Value = "0.07895"
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("D12").Value = "97.22"
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").Value = "1.866.82"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "5.151"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("D15").Value = "0.6809"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("D16").Value = "279.09"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("D17").Value = "30.211.33"
$ws.Range("E17").Value = "  +0.34%  "
$ws.Range("D18").Value = "13.66"
$ws.Range("E18").Value = "  +8.87%  "
$ws.Range("D19").Value = "0.9998"
$ws.Range("D20").Value = "0.000007321"
$ws.Range("E20").Value = "  +1.56%  "
$ws.Range("D21").Value = "2.112.24"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.370"
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("D23").Value = "0.9998"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "6.166"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").Value = "167.83"
$ws.Range("E25").Value = "  +0.58%  "
$ws.Range("D26").Value = "9.239"
$ws.Range("E26").Value = "  -0.83%  "
$ws.Range("D27").Value = "19.09"
$ws.Range("E27").Value = "  +1.14%  "
$ws.Range("D28").Value = "1.923"
$ws.Range("E28").Value = "  +0.90%  "
$ws.Range("E29").Value = "  +3.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09740"
$ws.Range("E30").Value = "  +1.97%  "
$ws.Range("D31").Value = "4.383"
$ws.Range("E31").Value = "  -0.95%  "
$ws.Range("D32").Value = "1.479"
$ws.Range("E32").Value = "  +0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.050"
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("D34").Value = "0.04734"
$ws.Range("E34").Value = "  +2.30%  "
$ws.Range("E35").Value = "  +3.25%  "
$ws.Range("D36").Value = "0.7086"
$ws.Range("E36").Value = "  +1.52%  "
$ws.Range("D37").Value = "2.707"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("E38").Value = "  +1.01%  "
$ws.Range("D39").Value = "2.606"
$ws.Range("E39").Value = "  +4.13%  "
$ws.Range("D40").Value = "6.279"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").Value = "75.16"
$ws.Range("E41").Value = "  +4.10%  "
$ws.Range("D42").Value = "1.953"
$ws.Range("E42").Value = "  +2.61%  "
$ws.Range("D43").Value = "0.8464"
$ws.Range("E43").Value = "  -0.81%  "
$ws.Range("D44").Value = "0.4176"
$ws.Range("E44").Value = "  +0.83%  "
$ws.Range("D45").Value = "0.9994"
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").Value = "103.19"
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("D47").Value = "963.87"
$ws.Range("E47").Value = "  -2.47%  "
$ws.Range("D48").Value = "7.205"
$ws.Range("E48").Value = "  +0.34%  "
$ws.Range("D49").Value = "9.329"
$ws.Range("E49").Value = "  +1.10%  "
$ws.Range("D50").Value = "34.13"
$ws.Range("E50").Value = "  +0.71%  "
$ws.Range("D51").Value = "0.05641"
$ws.Range("E51").Value = "  +0.41%  "
